$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 4.141153333333333
$ws.Cells.Item(2, 8).Value = 12.42346
$ws.Cells.Item(2, 9).Value = 0.2530231305454066
$ws.Cells.Item(2, 10).Value = 0.2530231305454066
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.100491
$ws.Cells.Item(2, 14).Value = 0.301473
$ws.Cells.Item(2, 15).Value = 0.03802271729753455
$ws.Cells.Item(2, 16).Value = 0.03802271729753454
$ws.Cells.Item(2, 17).Value = 0.4161486396199999
$ws.Cells.Item(2, 18).Value = 3.74533775658
$ws.Cells.Item(2, 19).Value = 0.009620626962465175
$ws.Cells.Item(2, 20).Value = 0.009620626962465174
$ws.Cells.Item(3, 7).Value = 4.141153333333333
$ws.Cells.Item(3, 8).Value = 12.42346
$ws.Cells.Item(3, 9).Value = 0.2530231305454066
$ws.Cells.Item(3, 10).Value = 0.2530231305454066
$ws.Cells.Item(3, 15).Value = 0.7272263254279357
$ws.Cells.Item(3, 16).Value = 0.7272263254279357
$ws.Cells.Item(3, 17).Value = 7.959300847819998
$ws.Cells.Item(3, 18).Value = 71.63370763037999
$ws.Cells.Item(3, 19).Value = 0.184005081474809
$ws.Cells.Item(3, 20).Value = 0.184005081474809
$ws.Cells.Item(4, 7).Value = 4.141153333333333
$ws.Cells.Item(4, 8).Value = 12.42346
$ws.Cells.Item(4, 9).Value = 0.2530231305454066
$ws.Cells.Item(4, 10).Value = 0.2530231305454066
$ws.Cells.Item(4, 13).Value = 0.6046453333333334
$ws.Cells.Item(4, 14).Value = 1.813936
$ws.Cells.Item(4, 15).Value = 0.2287792794838033
$ws.Cells.Item(4, 16).Value = 0.2287792794838033
$ws.Cells.Item(4, 17).Value = 2.503929037617778
$ws.Cells.Item(4, 18).Value = 22.53536133856
$ws.Cells.Item(4, 19).Value = 0.05788644949891444
$ws.Cells.Item(4, 20).Value = 0.05788644949891443
$ws.Cells.Item(5, 7).Value = 4.141153333333333
$ws.Cells.Item(5, 8).Value = 12.42346
$ws.Cells.Item(5, 9).Value = 0.2530231305454066
$ws.Cells.Item(5, 10).Value = 0.2530231305454066
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.01578266666666667
$ws.Cells.Item(5, 14).Value = 0.047348
$ws.Cells.Item(5, 15).Value = 0.005971677790726419
$ws.Cells.Item(5, 16).Value = 0.005971677790726419
$ws.Cells.Item(5, 17).Value = 0.06535844267555554
$ws.Cells.Item(5, 18).Value = 0.5882259840799999
$ws.Cells.Item(5, 19).Value = 0.001510972609218076
$ws.Cells.Item(5, 20).Value = 0.001510972609218076
$ws.Cells.Item(6, 9).Value = 0.3583796455306321
$ws.Cells.Item(6, 10).Value = 0.358379645530632
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.100491
$ws.Cells.Item(6, 14).Value = 0.301473
$ws.Cells.Item(6, 15).Value = 0.03802271729753455
$ws.Cells.Item(6, 16).Value = 0.03802271729753454
$ws.Cells.Item(6, 17).Value = 0.5894291230750001
$ws.Cells.Item(6, 18).Value = 5.304862107675
$ws.Cells.Item(6, 19).Value = 0.01362656794720186
$ws.Cells.Item(6, 20).Value = 0.01362656794720186
$ws.Cells.Item(7, 9).Value = 0.3583796455306321
$ws.Cells.Item(7, 10).Value = 0.358379645530632
$ws.Cells.Item(7, 15).Value = 0.7272263254279357
$ws.Cells.Item(7, 16).Value = 0.7272263254279357
$ws.Cells.Item(7, 19).Value = 0.2606231127274077
$ws.Cells.Item(7, 20).Value = 0.2606231127274076
$ws.Cells.Item(8, 9).Value = 0.3583796455306321
$ws.Cells.Item(8, 10).Value = 0.358379645530632
$ws.Cells.Item(8, 13).Value = 0.6046453333333334
$ws.Cells.Item(8, 14).Value = 1.813936
$ws.Cells.Item(8, 15).Value = 0.2287792794838033
$ws.Cells.Item(8, 16).Value = 0.2287792794838033
$ws.Cells.Item(8, 17).Value = 3.546542163955556
$ws.Cells.Item(8, 18).Value = 31.9188794756
$ws.Cells.Item(8, 19).Value = 0.08198983708615883
$ws.Cells.Item(8, 20).Value = 0.08198983708615881
$ws.Cells.Item(9, 9).Value = 0.3583796455306321
$ws.Cells.Item(9, 10).Value = 0.358379645530632
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.01578266666666667
$ws.Cells.Item(9, 14).Value = 0.047348
$ws.Cells.Item(9, 15).Value = 0.005971677790726419
$ws.Cells.Item(9, 16).Value = 0.005971677790726419
$ws.Cells.Item(9, 17).Value = 0.09257309981111113
$ws.Cells.Item(9, 18).Value = 0.8331578983000001
$ws.Cells.Item(9, 19).Value = 0.002140127769863682
$ws.Cells.Item(9, 20).Value = 0.002140127769863682
$ws.Cells.Item(10, 7).Value = 0.467591
$ws.Cells.Item(10, 8).Value = 1.402773
$ws.Cells.Item(10, 9).Value = 0.02856965900840602
$ws.Cells.Item(10, 10).Value = 0.02856965900840601
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.100491
$ws.Cells.Item(10, 14).Value = 0.301473
$ws.Cells.Item(10, 15).Value = 0.03802271729753455
$ws.Cells.Item(10, 16).Value = 0.03802271729753454
$ws.Cells.Item(10, 17).Value = 0.046988687181
$ws.Cells.Item(10, 18).Value = 0.422898184629
$ws.Cells.Item(10, 19).Value = 0.001086296067763583
$ws.Cells.Item(10, 20).Value = 0.001086296067763583
$ws.Cells.Item(11, 7).Value = 0.467591
$ws.Cells.Item(11, 8).Value = 1.402773
$ws.Cells.Item(11, 9).Value = 0.02856965900840602
$ws.Cells.Item(11, 10).Value = 0.02856965900840601
$ws.Cells.Item(11, 15).Value = 0.7272263254279357
$ws.Cells.Item(11, 16).Value = 0.7272263254279357
$ws.Cells.Item(11, 17).Value = 0.898710369591
$ws.Cells.Item(11, 18).Value = 8.088393326319
$ws.Cells.Item(11, 19).Value = 0.02077660813941223
$ws.Cells.Item(11, 20).Value = 0.02077660813941223
$ws.Cells.Item(12, 7).Value = 0.467591
$ws.Cells.Item(12, 8).Value = 1.402773
$ws.Cells.Item(12, 9).Value = 0.02856965900840602
$ws.Cells.Item(12, 10).Value = 0.02856965900840601
$ws.Cells.Item(12, 13).Value = 0.6046453333333334
$ws.Cells.Item(12, 14).Value = 1.813936
$ws.Cells.Item(12, 15).Value = 0.2287792794838033
$ws.Cells.Item(12, 16).Value = 0.2287792794838033
$ws.Cells.Item(12, 17).Value = 0.2827267160586667
$ws.Cells.Item(12, 18).Value = 2.544540444528
$ws.Cells.Item(12, 19).Value = 0.006536146003041079
$ws.Cells.Item(12, 20).Value = 0.006536146003041077
$ws.Cells.Item(13, 7).Value = 0.467591
$ws.Cells.Item(13, 8).Value = 1.402773
$ws.Cells.Item(13, 9).Value = 0.02856965900840602
$ws.Cells.Item(13, 10).Value = 0.02856965900840601
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.01578266666666667
$ws.Cells.Item(13, 14).Value = 0.047348
$ws.Cells.Item(13, 15).Value = 0.005971677790726419
$ws.Cells.Item(13, 16).Value = 0.005971677790726419
$ws.Cells.Item(13, 17).Value = 0.007379832889333334
$ws.Cells.Item(13, 18).Value = 0.066418496004
$ws.Cells.Item(13, 19).Value = 0.0001706087981891252
$ws.Cells.Item(13, 20).Value = 0.0001706087981891252
$ws.Cells.Item(14, 7).Value = 5.892462666666667
$ws.Cells.Item(14, 8).Value = 17.677388
$ws.Cells.Item(14, 9).Value = 0.3600275649155554
$ws.Cells.Item(14, 10).Value = 0.3600275649155554
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.100491
$ws.Cells.Item(14, 14).Value = 0.301473
$ws.Cells.Item(14, 15).Value = 0.03802271729753455
$ws.Cells.Item(14, 16).Value = 0.03802271729753454
$ws.Cells.Item(14, 17).Value = 0.592139465836
$ws.Cells.Item(14, 18).Value = 5.329255192524
$ws.Cells.Item(14, 19).Value = 0.01368922632010393
$ws.Cells.Item(14, 20).Value = 0.01368922632010393
$ws.Cells.Item(15, 7).Value = 5.892462666666667
$ws.Cells.Item(15, 8).Value = 17.677388
$ws.Cells.Item(15, 9).Value = 0.3600275649155554
$ws.Cells.Item(15, 10).Value = 0.3600275649155554
$ws.Cells.Item(15, 15).Value = 0.7272263254279357
$ws.Cells.Item(15, 16).Value = 0.7272263254279357
$ws.Cells.Item(15, 17).Value = 11.325319137796
$ws.Cells.Item(15, 18).Value = 101.927872240164
$ws.Cells.Item(15, 19).Value = 0.261821523086307
$ws.Cells.Item(15, 20).Value = 0.261821523086307
$ws.Cells.Item(16, 7).Value = 5.892462666666667
$ws.Cells.Item(16, 8).Value = 17.677388
$ws.Cells.Item(16, 9).Value = 0.3600275649155554
$ws.Cells.Item(16, 10).Value = 0.3600275649155554
$ws.Cells.Item(16, 13).Value = 0.6046453333333334
$ws.Cells.Item(16, 14).Value = 1.813936
$ws.Cells.Item(16, 15).Value = 0.2287792794838033
$ws.Cells.Item(16, 16).Value = 0.2287792794838033
$ws.Cells.Item(16, 17).Value = 3.562850053240889
$ws.Cells.Item(16, 18).Value = 32.065650479168
$ws.Cells.Item(16, 19).Value = 0.08236684689568899
$ws.Cells.Item(16, 20).Value = 0.08236684689568897
$ws.Cells.Item(17, 7).Value = 5.892462666666667
$ws.Cells.Item(17, 8).Value = 17.677388
$ws.Cells.Item(17, 9).Value = 0.3600275649155554
$ws.Cells.Item(17, 10).Value = 0.3600275649155554
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.01578266666666667
$ws.Cells.Item(17, 14).Value = 0.047348
$ws.Cells.Item(17, 15).Value = 0.005971677790726419
$ws.Cells.Item(17, 16).Value = 0.005971677790726419
$ws.Cells.Item(17, 17).Value = 0.09299877411377779
$ws.Cells.Item(17, 18).Value = 0.8369889670240001
$ws.Cells.Item(17, 19).Value = 0.002149968613455536
$ws.Cells.Item(17, 20).Value = 0.002149968613455536
